# "Fixing scripts for watch list"
#
# Row 34 (TestCase_E33) and Row 35 (TestCase_E34) on the "Test Cases" sheet
# get a second, related test scenario appended to their Jira-id and
# Description columns (joined with the sheet's "||" multi-value separator).
# Both rows also grow to a 30pt row height with wrapped text to
# accommodate the longer combined description, matching the formatting
# already used by similar multi-scenario rows (e.g. row 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 34 : TestCase_E33 -------------------------------------------------
$ws.Range("B34").Value = "OPQA-321`n||OPQA-621"
$ws.Range("C34").Value = "Verify that anyone can see the public watchlists of a user on user's profile page||Verify that user1 is able to see a watchlist on user2's profile page,  once user2's private watchlist is made to public."

# --- Row 35 : TestCase_E34 -------------------------------------------------
$ws.Range("B35").Value = "OPQA-329`n||OPQA-621"
$ws.Range("C35").Value = "Verify that no one can see the private watchlists of a user on user's profile page||Verify that user1 is not able to see a watchlist on user2's profile page,  once user2's public watchlist is reverted to private."

# Wrap the Jira-id column (Description already wraps) and grow the rows so
# the merged multi-line content is fully visible, same as row 38's style.
$ws.Range("B34").WrapText = $true
$ws.Range("B35").WrapText = $true
$ws.Rows.Item(34).RowHeight = 30
$ws.Rows.Item(35).RowHeight = 30

# The active selection on the sheet moved back to the top of the data.
$ws.Range("B2").Select()
